$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.495.78'
$ws.Range('E2').Value = '  +3.27%  '
$ws.Range('D3').Value = '2.993.99'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'562.70"
$ws.Range('E5').Value = '  +2.76%  '
$ws.Range('D6').Value = "'138.75"
$ws.Range('E6').Value = '  +7.01%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = "'0.518"
$ws.Range('E8').Value = '  +1.69%  '
$ws.Range('D9').Value = '2.983.18'
$ws.Range('E9').Value = '  +2.52%  '
$ws.Range('D10').Value = "'0.133"
$ws.Range('E10').Value = '  +5.58%  '
$ws.Range('D11').Value = "'5.28"
$ws.Range('E11').Value = '  +11.94%  '
$ws.Range('D12').Value = "'0.453"
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('D13').Value = "'0.0000230"
$ws.Range('E13').Value = '  +5.46%  '
$ws.Range('D14').Value = "'33.87"
$ws.Range('E14').Value = '  +3.82%  '
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '3.496.19'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').Value = "'7.18"
$ws.Range('E17').Value = '  +4.43%  '
$ws.Range('D18').Value = '2.994.55'
$ws.Range('E18').Value = '  +2.57%  '
$ws.Range('D19').Value = '59.511.32'
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('D20').Value = "'435.27"
$ws.Range('E20').Value = '  +4.64%  '
$ws.Range('D21').Value = "'13.57"
$ws.Range('E21').Value = '  +2.57%  '
$ws.Range('E22').Value = '  +4.22%  '
$ws.Range('E23').Value = '  +2.40%  '
$ws.Range('D24').Value = "'7.09"
$ws.Range('E24').Value = '  +1.88%  '
$ws.Range('D25').Value = "'80.20"
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  +11.35%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = "'2.54"
$ws.Range('E29').Value = '  +3.62%  '
$ws.Range('E30').Value = '  +5.76%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.107"
$ws.Range('E31').Value = '  +10.70%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = "'6.25"
$ws.Range('E32').Value = '  +5.17%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'25.77"
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('D34').Value = '0.0₃0776'
$ws.Range('E34').Value = '  +13.80%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').Value = "'0.989"
$ws.Range('E35').Value = '  +6.21%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = "'5.89"
$ws.Range('E36').Value = '  +4.02%  '
$ws.Range('D37').Value = "'2.09"
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').Value = "'48.99"
$ws.Range('E38').Value = '  +1.84%  '
$ws.Range('D39').Value = "'8.56"
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('D40').Value = "'2.78"
$ws.Range('E40').Value = '  +9.20%  '
$ws.Range('D41').Value = "'400.17"
$ws.Range('E41').Value = '  +7.68%  '
$ws.Range('D42').Value = "'0.0353"
$ws.Range('E42').Value = '  +2.96%  '
$ws.Range('D43').Value = '2.763.86'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('E45').Value = '  +7.55%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = "'122.93"
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').Value = "'34.34"
$ws.Range('E48').Value = '  +19.07%  '
$ws.Range('D49').Value = "'2.02"
$ws.Range('E49').Value = '  +4.00%  '
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('D51').Value = "'23.61"
$ws.Range('E51').Value = '  +3.70%  '
